# Apply crypto price/volume refresh per commit: "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.891.05'
$ws.Range("E2").Value = '  -2.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.619.43'
$ws.Range("E3").Value = '  -3.39%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.34'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.97'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.618.44'
$ws.Range("E9").Value = '  -3.36%  '
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("E12").Value = '  -1.37%  '
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.36'
$ws.Range("E14").Value = '  -3.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.095.81'
$ws.Range("E15").Value = '  -3.39%  '
$ws.Range("E16").Value = '  -3.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.711.16'
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.620.08'
$ws.Range("E18").Value = '  -3.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.97'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.00'
$ws.Range("E20").Value = '  +4.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '357.67'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("E23").Value = '  -5.90%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.90'
$ws.Range("E24").Value = '  +7.94%  '
$ws.Range("E25").Value = '  -6.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.86'
$ws.Range("E27").Value = '  -3.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.755.58'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  -3.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '550.28'
$ws.Range("E31").Value = '  -4.62%  '
$ws.Range("E32").Value = '  -3.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.36'
$ws.Range("E33").Value = '  -3.96%  '
$ws.Range("E34").Value = '  -4.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("E37").Value = '  -5.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.97'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.11'
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.366'
$ws.Range("E40").Value = '  -2.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.19'
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("E42").Value = '  -4.60%  '
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("E45").Value = '  -5.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.28'
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0296'
$ws.Range("E47").Value = '  -4.78%  '
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '151.33'
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  -3.21%  '
